# Edit script: add "2022-Q4" sheet data and update "总计" summary sheet.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: Create the new "2022-Q4" worksheet by duplicating the existing
# "2022-Q3" worksheet (tab index 2). This preserves all formatting,
# styles, column widths, etc. The duplicate is inserted immediately
# before the source sheet, i.e. it becomes the new tab index 2, and
# "2022-Q3" (and everything after it) shifts one position to the right.
# ---------------------------------------------------------------------
$srcQ3 = $wb.Worksheets.Item(2)
$srcQ3.Copy($srcQ3)
$newQ4 = $wb.Worksheets.Item(2)
$newQ4.Name = "2022-Q4"

# The duplicated sheet has 7 rows (header + 6 funds); the new quarter
# only has 4 funds, so remove the two extra data rows (rows 6 and 7).
$newQ4.Rows.Item(7).Delete()
$newQ4.Rows.Item(6).Delete()

# Force text columns to stay text (preserve values such as leading
# zeros in fund codes, and numeric-looking percentages as strings)
$newQ4.Range("B2:B5").NumberFormat = "@"
$newQ4.Range("C2:C5").NumberFormat = "@"
$newQ4.Range("D2:G5").NumberFormat = "@"

# Row 2: 006679
$newQ4.Cells.Item(2,1).Value = 0
$newQ4.Cells.Item(2,2).Value = "006679"
$newQ4.Cells.Item(2,3).Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 A"
$newQ4.Cells.Item(2,4).Value = "11.33"
$newQ4.Cells.Item(2,5).Value = "90.17"
$newQ4.Cells.Item(2,6).Value = "8.68"
$newQ4.Cells.Item(2,7).Value = "0.9834"
$newQ4.Cells.Item(2,8).Value = 2

# Row 3: 162719
$newQ4.Cells.Item(3,1).Value = 1
$newQ4.Cells.Item(3,2).Value = "162719"
$newQ4.Cells.Item(3,3).Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A"
$newQ4.Cells.Item(3,4).Value = "11.33"
$newQ4.Cells.Item(3,5).Value = "90.17"
$newQ4.Cells.Item(3,6).Value = "8.68"
$newQ4.Cells.Item(3,7).Value = "0.9834"
$newQ4.Cells.Item(3,8).Value = 2

# Row 4: 006680
$newQ4.Cells.Item(4,1).Value = 2
$newQ4.Cells.Item(4,2).Value = "006680"
$newQ4.Cells.Item(4,3).Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 C"
$newQ4.Cells.Item(4,4).Value = "5.56"
$newQ4.Cells.Item(4,5).Value = "90.17"
$newQ4.Cells.Item(4,6).Value = "8.68"
$newQ4.Cells.Item(4,7).Value = "0.4826"
$newQ4.Cells.Item(4,8).Value = 2

# Row 5: 004243
$newQ4.Cells.Item(5,1).Value = 3
$newQ4.Cells.Item(5,2).Value = "004243"
$newQ4.Cells.Item(5,3).Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C"
$newQ4.Cells.Item(5,4).Value = "-11.33"
$newQ4.Cells.Item(5,5).Value = "90.17"
$newQ4.Cells.Item(5,6).Value = "8.68"
$newQ4.Cells.Item(5,7).Value = "-0.9834"
$newQ4.Cells.Item(5,8).Value = 2

# ---------------------------------------------------------------------
# Step 2: Update the "总计" (summary) sheet - insert a new row for
# "2022-Q4" above the existing "2022-Q3" row, pushing all other
# quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# Clear any carried-over formatting on the data cells, then restore the
# bold/bordered look used by the rest of column A.
$summary.Range("B2:D2").ClearFormats()
$aCell = $summary.Cells.Item(2,1)
$aCell.Font.Bold = $true
$aCell.Borders.LineStyle = 1
$aCell.HorizontalAlignment = -4108
$aCell.VerticalAlignment = -4160

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 4
$summary.Cells.Item(2,4).Value = 1.47

# ---------------------------------------------------------------------
# Step 3: Keep "2020-Q4" as the active/selected sheet (it was the
# active sheet before the edit, and its tab index shifts from 9 to 10).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
